# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period detail table (rows 16-27, columns C:E) is regrouped so
# that all six workers appear first for period 1805, then all six workers
# again for period 1806 (previously the two periods were interleaved per
# worker). Column C = worker document number, D = worker name, E = period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of workers as they appear in the table (doc number, name).
$workers = @(
    @("1047417915", "CARLOS ALFONSO PALENCIA RODRIGUEZ"),
    @("1143363639", "OSCAR ALFONSO PALENCIA RODRIGUEZ"),
    @("1047431310", "JORGE ELIECER PALENCIA RODRIGUEZ"),
    @("1143353062", "HUGUER ENRIQUE PALENCIA RODRIGUEZ"),
    @("1047416352", "JHON JAIRO PALENCIA RODRIGUEZ"),
    @("19935236",   "EUGENIO VALDES MOSQUERA")
)

$periods = @("1805", "1806")

$row = 16
foreach ($period in $periods) {
    foreach ($worker in $workers) {
        $ws.Range("C$row").Value = $worker[0]
        $ws.Range("D$row").Value = $worker[1]
        $ws.Range("E$row").Value = $period
        $row = $row + 1
    }
}
